# New weekly Pomelo price observation (Vega Central Mapocho de Santiago).
# Insert a new row at position 20 (pushing existing rows 20..54 down to
# 21..55) and populate it with the new day's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value  = 9
$ws.Cells.Item(20, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value  = "Metropolitana"
$ws.Cells.Item(20, 4).Value  = 44482
$ws.Cells.Item(20, 5).Value  = 13
$ws.Cells.Item(20, 6).Value  = "Fruta"
$ws.Cells.Item(20, 7).Value  = 100102
$ws.Cells.Item(20, 8).Value  = "Cítricos"
$ws.Cells.Item(20, 9).Value  = 100102006
$ws.Cells.Item(20, 10).Value = "Pomelo"
$ws.Cells.Item(20, 11).Value = "Start Ruby"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 330
$ws.Cells.Item(20, 14).Value = 7000
$ws.Cells.Item(20, 15).Value = 7500
$ws.Cells.Item(20, 16).Value = 7273
$ws.Cells.Item(20, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(20, 18).Value = "Región Metropolitana"
$ws.Cells.Item(20, 19).Value = 520
$ws.Cells.Item(20, 20).Value = 14
